# Apply the change described by the diff: add columns I and J
# (headers "I0"/"IF" plus numeric data) to the existing sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from an existing header cell (H1) onto the new header cells
# so I1/J1 pick up the same bold/centered/bordered formatting (style index 1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
